$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows (rows 2-13), columns A-J
$data = @(
    @("male",   "June26", "peakValue",    416,    0.02284220143734341,   "mannwhitneyu", 0.02251127796471279,  "Dunn", 0.3301127214170693,   -30.77626140633384),
    @("male",   "June26", "tau",          853,    0.0099800429795172,    "mannwhitneyu", 0.009820543828851872, "Dunn", -0.3735909822866346,  28.14803918411168),
    @("male",   "May20",  "peakValue",    478,    0.02168335053511856,   "mannwhitneyu", 0.02139131675050924,  "Dunn", 0.3190883190883191,   -19.02877492877494),
    @("female", "April16","reactionTime", 914.5,  0.003309579419531099,  "mannwhitneyu", 0.003251122413952952, "Dunn", -0.4112654320987654,  0.7666666666666693),
    @("female", "April16","peakTime",     941.5,  0.0003618647162627166, "mannwhitneyu", 0.0003535553487708839,"Dunn", -0.4529320987654322,  0.9240740740740776),
    @("female", "April16","peakValue",    463,    0.04576783585618472,   "mannwhitneyu", 0.04518356364957534,  "Dunn", 0.2854938271604939,   -7.782407407407412),
    @("female", "June26", "peakValue",    391,    0.00280302378995304,   "mannwhitneyu", 0.002755054749030001, "Dunn", 0.4207407407407407,   -19.14172839506172),
    @("female", "June26", "tau",          1068,   0.00003513962043189032, "mannwhitneyu", 0.00003434149222123235,"Dunn", -0.5822222222222222,  15.91196049382717),
    @("female", "May20",  "peakTime",     883.5,  0.007059101550412708,  "mannwhitneyu", 0.006938634220435131, "Dunn", -0.3634259259259258,  -0.2962962962962976),
    @("female", "May20",  "peakValue",    397.5,  0.006795668520630508,  "mannwhitneyu", 0.006685675302743306, "Dunn", 0.3865740740740741,   -22.67160493827161),
    @("female", "May20",  "RMS",          379,    0.003650961021094545,  "mannwhitneyu", 0.003588280051067788, "Dunn", 0.4151234567901234,   -6.162679012345681),
    @("female", "May20",  "AUC",          359,    0.001788061142657725,  "mannwhitneyu", 0.001755449013430172, "Dunn", 0.4459876543209876,   -2480.932794753087)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
    $ws.Cells.Item($row, 9).Value = $vals[8]
    $ws.Cells.Item($row, 10).Value = $vals[9]
}
